$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 (China)
$ws.Cells.Item(2, 3).Value2 = 342
$ws.Cells.Item(2, 6).Value2 = 45444
$ws.Cells.Item(2, 7).Value2 = 30319
$ws.Cells.Item(2, 8).Value2 = 45505

# Row 3 (United States)
$ws.Cells.Item(3, 5).Value2 = 30256
$ws.Cells.Item(3, 6).Value2 = 45444

# Row 4 (Euro Area)
$ws.Cells.Item(4, 5).Value2 = 30256
$ws.Cells.Item(4, 6).Value2 = 45444
$ws.Cells.Item(4, 7).Value2 = 30319
$ws.Cells.Item(4, 8).Value2 = 45505

# Row 5 (Japan)
$ws.Cells.Item(5, 5).Value2 = 30256
$ws.Cells.Item(5, 6).Value2 = 45444
$ws.Cells.Item(5, 7).Value2 = 30319
$ws.Cells.Item(5, 8).Value2 = 45505

# Row 6 (United Kingdom)
$ws.Cells.Item(6, 3).Value2 = 451
$ws.Cells.Item(6, 6).Value2 = 45444
$ws.Cells.Item(6, 7).Value2 = 30319
$ws.Cells.Item(6, 8).Value2 = 45505

# Row 7 (South Korea)
$ws.Cells.Item(7, 5).Value2 = 30256
$ws.Cells.Item(7, 6).Value2 = 45444
$ws.Cells.Item(7, 7).Value2 = 30319
$ws.Cells.Item(7, 8).Value2 = 45505

# Row 8 (Hong Kong)
$ws.Cells.Item(8, 4).Value2 = 418
$ws.Cells.Item(8, 5).Value2 = 30256
$ws.Cells.Item(8, 6).Value2 = 45444
$ws.Cells.Item(8, 8).Value2 = 45505

# Row 9 (Australia)
$ws.Cells.Item(9, 5).Value2 = 30256
$ws.Cells.Item(9, 6).Value2 = 45444
$ws.Cells.Item(9, 7).Value2 = 30319
$ws.Cells.Item(9, 8).Value2 = 45505

# Row 10 (Taiwan)
$ws.Cells.Item(10, 4).Value2 = 491
$ws.Cells.Item(10, 5).Value2 = 30256
$ws.Cells.Item(10, 6).Value2 = 45444
$ws.Cells.Item(10, 8).Value2 = 45505

# Row 11 (Canada)
$ws.Cells.Item(11, 5).Value2 = 30225
$ws.Cells.Item(11, 6).Value2 = 45413
$ws.Cells.Item(11, 7).Value2 = 30319
$ws.Cells.Item(11, 8).Value2 = 45505

# Row 12 (Russia)
$ws.Cells.Item(12, 3).Value2 = 379
$ws.Cells.Item(12, 4).Value2 = 360
$ws.Cells.Item(12, 6).Value2 = 45444
$ws.Cells.Item(12, 8).Value2 = 45505

# Row 13 (Switzerland)
$ws.Cells.Item(13, 3).Value2 = 474
$ws.Cells.Item(13, 6).Value2 = 45413
$ws.Cells.Item(13, 7).Value2 = 30319
$ws.Cells.Item(13, 8).Value2 = 45505

# Row 14 (Brazil)
$ws.Cells.Item(14, 4).Value2 = 404
$ws.Cells.Item(14, 8).Value2 = 45505

# Row 15 (India)
$ws.Cells.Item(15, 7).Value2 = 30286
$ws.Cells.Item(15, 8).Value2 = 45505

# Row 16 (Mexico)
$ws.Cells.Item(16, 3).Value2 = 462
$ws.Cells.Item(16, 4).Value2 = 418
$ws.Cells.Item(16, 6).Value2 = 45413
$ws.Cells.Item(16, 8).Value2 = 45505

# Row 17 (Saudi Arabia)
$ws.Cells.Item(17, 3).Value2 = 378
$ws.Cells.Item(17, 4).Value2 = 402
$ws.Cells.Item(17, 6).Value2 = 45444
$ws.Cells.Item(17, 8).Value2 = 45505

# Row 18 (Singapore)
$ws.Cells.Item(18, 4).Value2 = 268
$ws.Cells.Item(18, 5).Value2 = 30256
$ws.Cells.Item(18, 6).Value2 = 45444
$ws.Cells.Item(18, 8).Value2 = 45505

# Row 19 (Indonesia)
$ws.Cells.Item(19, 4).Value2 = 406
$ws.Cells.Item(19, 5).Value2 = 30256
$ws.Cells.Item(19, 6).Value2 = 45444
$ws.Cells.Item(19, 8).Value2 = 45505

# Row 20 (Malaysia)
$ws.Cells.Item(20, 3).Value2 = 489
$ws.Cells.Item(20, 6).Value2 = 45444
$ws.Cells.Item(20, 7).Value2 = 30319
$ws.Cells.Item(20, 8).Value2 = 45505

# Row 21 (Norway)
$ws.Cells.Item(21, 5).Value2 = 30225
$ws.Cells.Item(21, 6).Value2 = 45413
$ws.Cells.Item(21, 7).Value2 = 30319
$ws.Cells.Item(21, 8).Value2 = 45505

# Row 22 (Philippines)
$ws.Cells.Item(22, 4).Value2 = 388
$ws.Cells.Item(22, 5).Value2 = 30256
$ws.Cells.Item(22, 6).Value2 = 45444
$ws.Cells.Item(22, 8).Value2 = 45505

# Row 23 (New Zealand)
$ws.Cells.Item(23, 4).Value2 = 225
$ws.Cells.Item(23, 5).Value2 = 30225
$ws.Cells.Item(23, 6).Value2 = 45413
$ws.Cells.Item(23, 8).Value2 = 45505

# Row 24 (Denmark)
$ws.Cells.Item(24, 3).Value2 = 402
$ws.Cells.Item(24, 6).Value2 = 45444
$ws.Cells.Item(24, 7).Value2 = 30319
$ws.Cells.Item(24, 8).Value2 = 45505

# Row 25 (South Africa)
$ws.Cells.Item(25, 7).Value2 = 30319
$ws.Cells.Item(25, 8).Value2 = 45505

# Row 26 (Chile)
$ws.Cells.Item(26, 3).Value2 = 462
$ws.Cells.Item(26, 4).Value2 = 406
$ws.Cells.Item(26, 6).Value2 = 45444
$ws.Cells.Item(26, 8).Value2 = 45505

# Row 27 (Colombia)
$ws.Cells.Item(27, 4).Value2 = 418
$ws.Cells.Item(27, 5).Value2 = 30256
$ws.Cells.Item(27, 6).Value2 = 45444
$ws.Cells.Item(27, 8).Value2 = 45505

# Row 28 (Kuwait)
$ws.Cells.Item(28, 3).Value2 = 367
$ws.Cells.Item(28, 4).Value2 = 375
$ws.Cells.Item(28, 6).Value2 = 45444
$ws.Cells.Item(28, 8).Value2 = 45505
